$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Главные")

$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 637
$ws.Range("E2").Value = 267
$ws.Range("F2").Value = 370
$ws.Range("G2").Value = 21.23
$ws.Range("H2").Value = 8.9
$ws.Range("I2").Value = 12.33
$ws.Range("J2").Value = 116
$ws.Range("K2").Value = 140
$ws.Range("V2").Value = 10
$ws.Range("AA2").Value = "2025-11-27 03:06:04"

$ws.Range("C3").Value = 28
$ws.Range("D3").Value = 501
$ws.Range("E3").Value = 229
$ws.Range("F3").Value = 272
$ws.Range("G3").Value = 17.89
$ws.Range("H3").Value = 8.18
$ws.Range("I3").Value = 9.710000000000001
$ws.Range("J3").Value = 112
$ws.Range("K3").Value = 111
$ws.Range("M3").Value = 4
$ws.Range("W3").Value = 10
$ws.Range("Y3").Value = 4
$ws.Range("AA3").Value = "2025-11-27 03:06:04"

$ws.Range("AA4").Value = "2025-11-27 03:06:04"

$ws.Range("AA5").Value = "2025-11-27 03:06:04"

$ws.Range("C6").Value = 29
$ws.Range("D6").Value = 483
$ws.Range("E6").Value = 215
$ws.Range("F6").Value = 268
$ws.Range("G6").Value = 16.66
$ws.Range("H6").Value = 7.41
$ws.Range("I6").Value = 9.24
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 114
$ws.Range("V6").Value = 8
$ws.Range("AA6").Value = "2025-11-27 03:06:04"

$ws.Range("AA7").Value = "2025-11-27 03:06:04"

$ws.Range("AA8").Value = "2025-11-27 03:06:04"

$ws.Range("AA9").Value = "2025-11-27 03:06:04"

$ws.Range("AA10").Value = "2025-11-27 03:06:04"

$ws.Range("C11").Value = 21
$ws.Range("D11").Value = 516
$ws.Range("E11").Value = 234
$ws.Range("F11").Value = 282
$ws.Range("G11").Value = 24.57
$ws.Range("H11").Value = 11.14
$ws.Range("I11").Value = 13.43
$ws.Range("J11").Value = 102
$ws.Range("K11").Value = 96
$ws.Range("AA11").Value = "2025-11-27 03:06:04"

$ws.Range("AA12").Value = "2025-11-27 03:06:04"

$ws.Range("AA13").Value = "2025-11-27 03:06:04"

$ws.Range("C14").Value = 20
$ws.Range("D14").Value = 268
$ws.Range("E14").Value = 140
$ws.Range("F14").Value = 128
$ws.Range("G14").Value = 13.4
$ws.Range("H14").Value = 7
$ws.Range("I14").Value = 6.4
$ws.Range("J14").Value = 65
$ws.Range("K14").Value = 54
$ws.Range("L14").Value = 2
$ws.Range("M14").Value = 4
$ws.Range("W14").Value = 14
$ws.Range("AA14").Value = "2025-11-27 03:06:04"

$ws.Range("AA15").Value = "2025-11-27 03:06:04"

$ws.Range("AA16").Value = "2025-11-27 03:06:04"

$ws.Range("AA17").Value = "2025-11-27 03:06:04"

$ws.Range("AA18").Value = "2025-11-27 03:06:04"

$ws.Range("AA19").Value = "2025-11-27 03:06:04"

$ws.Range("AA20").Value = "2025-11-27 03:06:04"

$ws.Range("AA21").Value = "2025-11-27 03:06:04"

$ws.Range("AA22").Value = "2025-11-27 03:06:04"

$ws.Range("AA23").Value = "2025-11-27 03:06:04"

$ws.Range("AA24").Value = "2025-11-27 03:06:04"

$ws.Range("AA25").Value = "2025-11-27 03:06:04"

$ws.Range("AA26").Value = "2025-11-27 03:06:04"


$ws = $wb.Worksheets.Item("Линейные")

$ws.Range("C2").Value = 18
$ws.Range("D2").Value = 338
$ws.Range("E2").Value = 150
$ws.Range("F2").Value = 188
$ws.Range("G2").Value = 18.78
$ws.Range("H2").Value = 8.33
$ws.Range("I2").Value = 10.44
$ws.Range("J2").Value = 65
$ws.Range("K2").Value = 69
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = 2
$ws.Range("W2").Value = 12
$ws.Range("AA2").Value = "2025-11-27 03:06:04"

$ws.Range("AA3").Value = "2025-11-27 03:06:04"

$ws.Range("AA4").Value = "2025-11-27 03:06:04"

$ws.Range("AA5").Value = "2025-11-27 03:06:04"

$ws.Range("AA6").Value = "2025-11-27 03:06:04"

$ws.Range("AA7").Value = "2025-11-27 03:06:04"

$ws.Range("AA8").Value = "2025-11-27 03:06:04"

$ws.Range("AA9").Value = "2025-11-27 03:06:04"

$ws.Range("AA10").Value = "2025-11-27 03:06:04"

$ws.Range("C11").Value = 18
$ws.Range("D11").Value = 258
$ws.Range("E11").Value = 112
$ws.Range("F11").Value = 146
$ws.Range("G11").Value = 14.33
$ws.Range("H11").Value = 6.22
$ws.Range("I11").Value = 8.109999999999999
$ws.Range("J11").Value = 56
$ws.Range("K11").Value = 68
$ws.Range("AA11").Value = "2025-11-27 03:06:04"

$ws.Range("C12").Value = 23
$ws.Range("D12").Value = 402
$ws.Range("E12").Value = 191
$ws.Range("F12").Value = 211
$ws.Range("G12").Value = 17.48
$ws.Range("H12").Value = 8.300000000000001
$ws.Range("I12").Value = 9.17
$ws.Range("J12").Value = 88
$ws.Range("K12").Value = 98
$ws.Range("V12").Value = 12
$ws.Range("AA12").Value = "2025-11-27 03:06:04"

$ws.Range("AA13").Value = "2025-11-27 03:06:04"

$ws.Range("AA14").Value = "2025-11-27 03:06:04"

$ws.Range("AA15").Value = "2025-11-27 03:06:04"

$ws.Range("AA16").Value = "2025-11-27 03:06:04"

$ws.Range("AA17").Value = "2025-11-27 03:06:04"

$ws.Range("C18").Value = 30
$ws.Range("D18").Value = 499
$ws.Range("E18").Value = 232
$ws.Range("F18").Value = 267
$ws.Range("G18").Value = 16.63
$ws.Range("H18").Value = 7.73
$ws.Range("I18").Value = 8.9
$ws.Range("J18").Value = 111
$ws.Range("K18").Value = 116
$ws.Range("AA18").Value = "2025-11-27 03:06:04"

$ws.Range("AA19").Value = "2025-11-27 03:06:04"

$ws.Range("AA20").Value = "2025-11-27 03:06:04"

$ws.Range("AA21").Value = "2025-11-27 03:06:04"

$ws.Range("AA22").Value = "2025-11-27 03:06:04"

$ws.Range("C23").Value = 16
$ws.Range("D23").Value = 244
$ws.Range("E23").Value = 120
$ws.Range("F23").Value = 124
$ws.Range("G23").Value = 15.25
$ws.Range("H23").Value = 7.5
$ws.Range("I23").Value = 7.75
$ws.Range("J23").Value = 55
$ws.Range("K23").Value = 57
$ws.Range("L23").Value = 2
$ws.Range("M23").Value = 2
$ws.Range("W23").Value = 4
$ws.Range("AA23").Value = "2025-11-27 03:06:04"

$ws.Range("AA24").Value = "2025-11-27 03:06:04"

$ws.Range("AA25").Value = "2025-11-27 03:06:04"

$ws.Range("AA26").Value = "2025-11-27 03:06:04"

